$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Silver Rear_side")
$ws1.Range("B18").Value = "'5,454"

$ws2 = $wb.Worksheets.Item("Silver Busbar front-side")
$ws2.Range("B18").Value = "'8,166"

$ws3 = $wb.Worksheets.Item("Silver finger front-side")
$ws3.Range("B18").Value = "'8,216"

$ws4 = $wb.Worksheets.Item("USD_CNY")
$ws4.Range("B18").Value = "'7.2456"
